# Populate the example row (row 2) of the revenue template, and
# apply the red "Paid" highlight style to the spacer row (row 3).
#
# Shared-string insertion order matters for a byte-exact OOXML match,
# so the new strings are written in the same order Excel produced them:
# F2 ("Example") first, then A2, B2, E2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "Example"
$ws.Range("A2").Value = "G01 - Rent 12/2002"
$ws.Range("B2").Value = "G01 - Lease"
$ws.Range("C2").Value = 1000
$ws.Range("D2").Value = Get-Date -Year 2020 -Month 11 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Range("E2").Value = "Paid"

$ws.Rows.Item(2).RowHeight = 30

# Row 3 becomes a red-filled "spacer" row under the new example.
$ws.Range("A3:E3").Interior.Color = 255
$ws.Range("A3:F3").Select() | Out-Null
$ws.Range("A3").Select() | Out-Null
